$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Octubre de 2020 a las 16:27"

# Row data: country name (col A) plus Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes (cols B-H)
$data = @(
    @(4, "Estados Unidos", 8524139, 3189, 5548054, 2749812, 0, 89, 226273),
    @(5, "India", 7669023, 19865, 6811403, 741568, 0, 102, 116052),
    @(21, "Alemania", 383016, 2118, 302100, 70954, 0, 7, 9962),
    @(39, "Catar", 130210, 266, 127093, 2892, 0, 1, 225),
    @(47, "Portugal", 106271, 2535, 63238, 40804, 0, 16, 2229),
    @(48, "Egipto", 105705, 0, 98413, 1150, 0, 0, 6142),
    @(52, "Suiza", 91763, 5596, 54600, 35136, 0, 5, 2027),
    @(55, "Bielorrusia", 89642, 733, 80905, 7796, 0, 4, 941),
    @(65, "Singapur", 57933, 12, 57821, 84, 0, 0, 28),
    @(75, "Kenia", 46144, 497, 32760, 12526, 0, 16, 858),
    @(79, "Birmania", 39696, 1194, 18865, 19859, 0, 27, 972),
    @(96, "Noruega", 16880, 109, 11863, 4738, 0, 1, 279),
    @(97, "Madagascar", 16810, 0, 16215, 357, 0, 0, 238),
    @(105, "Guinea", 11599, 61, 10461, 1068, 0, 0, 70),
    @(111, "Tayikistan", 10613, 39, 9668, 865, 0, 0, 80),
    @(182, "Islas Feroe", 488, 0, 473, 15, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = [double]$row[2]
    $ws.Cells.Item($r, 3).Value = [double]$row[3]
    $ws.Cells.Item($r, 4).Value = [double]$row[4]
    $ws.Cells.Item($r, 5).Value = [double]$row[5]
    $ws.Cells.Item($r, 6).Value = [double]$row[6]
    $ws.Cells.Item($r, 7).Value = [double]$row[7]
    $ws.Cells.Item($r, 8).Value = [double]$row[8]
}
